$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 525, pushing the existing weekly
# records (old rows 525-551) down to rows 526-552.
$ws.Rows.Item(525).Insert()

# Populate the newly inserted row 525 with this week's price record.
$ws.Cells.Item(525, 1).Value = 3
$ws.Cells.Item(525, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(525, 3).Value = 'Coquimbo'
$ws.Cells.Item(525, 4).Value = 45041
$ws.Cells.Item(525, 5).Value = 5
$ws.Cells.Item(525, 6).Value = 100112043
$ws.Cells.Item(525, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(525, 8).Value = 'Sin especificar'
$ws.Cells.Item(525, 9).Value = 'Primera'
$ws.Cells.Item(525, 10).Value = 65
$ws.Cells.Item(525, 11).Value = 9000
$ws.Cells.Item(525, 12).Value = 9000
$ws.Cells.Item(525, 13).Value = 9000
$ws.Cells.Item(525, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(525, 15).Value = 'Limache'
$ws.Cells.Item(525, 16).Value = 150
$ws.Cells.Item(525, 17).Value = 60
$ws.Cells.Item(525, 18).Value = 'Hortaliza'
